# Apply parameter-boundary updates to the LacI model "parameters_0" workbook.
# Column B = bmin, Column C = bmax. Originally several C cells held a
# formula "=B*1000" (partly as a shared formula group); after manual
# parameter estimation the boundaries were narrowed down and entered as
# independent literal values, so every C-cell touched here becomes a
# plain numeric value (this also breaks/removes the old formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (bmin) updates ---
$ws.Range("B2").Value  = 0.1
$ws.Range("B3").Value  = 15
$ws.Range("B4").Value  = 0.0001
$ws.Range("B5").Value  = 0.0001
$ws.Range("B6").Value  = 0.01
$ws.Range("B7").Value  = 1
$ws.Range("B11").Value = 0.001
$ws.Range("B15").Value = 1
$ws.Range("B17").Value = 0.00001
$ws.Range("B19").Value = 1
$ws.Range("B20").Value = 5
$ws.Range("B21").Value = 0.001
$ws.Range("B22").Value = 1

# --- Column C (bmax) updates: now plain values (no more =B*1000 formula) ---
$ws.Range("C2").Value  = 3
$ws.Range("C3").Value  = 25
$ws.Range("C4").Value  = 0.1
$ws.Range("C5").Value  = 0.1
$ws.Range("C6").Value  = 0.05
$ws.Range("C7").Value  = 100
$ws.Range("C8").Value  = 5
$ws.Range("C9").Value  = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 100
$ws.Range("C12").Value = 0.01
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 10
$ws.Range("C16").Value = 0.01
$ws.Range("C17").Value = 0.001
$ws.Range("C19").Value = 10
$ws.Range("C20").Value = 25
$ws.Range("C21").Value = 0.01
$ws.Range("C22").Value = 10
$ws.Range("C23").Value = 3

# --- Active selection moved from B22 to C8 ---
$ws.Range("C8").Select()
